# Update edited session - rename sheet and remove a stray duplicate scan row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Scanner" to "Session"
$ws.Name = "Session"

# Remove the original row 33 (Student ID 244484 / 11:18:36 Scan) which was
# a duplicate/erroneous scan entry. This shifts rows 34:64 up to 33:63 and
# reduces the used range to A1:F63.
$ws.Rows.Item(33).Delete()
